$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.353242993354797
$ws.Range("B1").Value = 3.427502155303955
$ws.Range("C1").Value = 5.281919479370117
$ws.Range("D1").Value = 7.505433082580566
$ws.Range("E1").Value = 2.738493204116821
